$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.222.25"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").Value = "1.683.47"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("D5").Value = "'215.63"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").Value = "'0.519"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "'0.257"
$ws.Range("E8").Value = "  +2.13%  "
$ws.Range("D9").Value = "'21.65"
$ws.Range("E9").Value = "  +6.42%  "
$ws.Range("E10").Value = "  +0.56%  "
$ws.Range("D11").Value = "'0.0890"
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").Value = "1.919.71"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").Value = "1.682.31"
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("E14").Value = "  +1.72%  "
$ws.Range("D15").Value = "'0.543"
$ws.Range("E15").Value = "  +2.76%  "
$ws.Range("D16").Value = "'66.44"
$ws.Range("E16").Value = "  +0.85%  "
$ws.Range("D17").Value = "27.198.08"
$ws.Range("E17").Value = "  +0.78%  "
$ws.Range("D18").Value = "'239.31"
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("D19").Value = "'8.09"
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("E20").Value = "  +1.52%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").Value = "'4.56"
$ws.Range("E22").Value = "  +2.57%  "
$ws.Range("D23").Value = "'9.50"
$ws.Range("E23").Value = "  +3.19%  "
$ws.Range("E24").Value = "  -3.50%  "
$ws.Range("D25").Value = "'148.11"
$ws.Range("E25").Value = "  +1.74%  "
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("D27").Value = "'16.38"
$ws.Range("E27").Value = "  +1.88%  "
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("E31").Value = "  +0.37%  "
$ws.Range("D32").Value = "1.573.53"
$ws.Range("E32").Value = "  +5.75%  "
$ws.Range("E33").Value = "  +1.52%  "
$ws.Range("E34").Value = "  +2.62%  "
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("D36").Value = "'0.603"
$ws.Range("E36").Value = "  +2.77%  "
$ws.Range("D37").Value = "'0.943"
$ws.Range("E37").Value = "  +4.59%  "
$ws.Range("E38").Value = "  -1.05%  "
$ws.Range("E39").Value = "  -0.33%  "
$ws.Range("E40").Value = "  +3.89%  "
$ws.Range("D41").Value = "'69.18"
$ws.Range("E41").Value = "  +2.33%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  -4.23%  "
$ws.Range("E44").Value = "  -2.33%  "
$ws.Range("D45").Value = "1.829.09"
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("E46").Value = "  +0.91%  "
$ws.Range("D47").Value = "'90.96"
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("E49").Value = "  +1.88%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'8.19"
$ws.Range("E50").Value = "  +6.15%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.104"
$ws.Range("E51").Value = "  +1.91%  "
